$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/27/2023  Through  12/3/2023"

# --- Cells changing from text placeholder to numeric (apply matching number format first) ---
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 2
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = 0
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = -50

# --- Cells changing from numeric to text placeholder ---
# Force text type via "@" format, assign, then copy the exact visual style (s=14) from a known text cell
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("C28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null

# --- Plain numeric value updates (no type/style change) ---
$ws.Range("M14").Value = -20
$ws.Range("I15").Value = 33
$ws.Range("J15").Value = 36
$ws.Range("K15").Value = -8.333333333333
$ws.Range("L15").Value = 17.857142857142
$ws.Range("M15").Value = 26.923076923076
$ws.Range("N15").Value = -43.103448275862
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 30
$ws.Range("H16").Value = 57.894736842105
$ws.Range("I16").Value = 426
$ws.Range("J16").Value = 368
$ws.Range("K16").Value = 15.760869565217
$ws.Range("L16").Value = 25.294117647058
$ws.Range("M16").Value = -9.168443496801
$ws.Range("N16").Value = -73.207547169811
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 63.636363636363
$ws.Range("F17").Value = 54
$ws.Range("H17").Value = 17.391304347826
$ws.Range("I17").Value = 643
$ws.Range("J17").Value = 567
$ws.Range("K17").Value = 13.403880070546
$ws.Range("L17").Value = 26.824457593688
$ws.Range("M17").Value = 51.650943396226
$ws.Range("N17").Value = 7.885906040268
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 20
$ws.Range("H18").Value = 53.846153846153
$ws.Range("I18").Value = 212
$ws.Range("J18").Value = 178
$ws.Range("K18").Value = 19.101123595505
$ws.Range("L18").Value = 29.268292682926
$ws.Range("M18").Value = -45.780051150895
$ws.Range("N18").Value = -90.315212425765
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 157.142857142857
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 22.641509433962
$ws.Range("I19").Value = 678
$ws.Range("J19").Value = 614
$ws.Range("K19").Value = 10.423452768729
$ws.Range("L19").Value = 9.354838709677
$ws.Range("M19").Value = 29.63671128107
$ws.Range("N19").Value = -20.04716981132
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -14.285714285714
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -5
$ws.Range("I20").Value = 313
$ws.Range("J20").Value = 238
$ws.Range("K20").Value = 31.512605042016
$ws.Range("L20").Value = 54.950495049505
$ws.Range("M20").Value = 86.309523809523
$ws.Range("N20").Value = -77.335264301231
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 77.419354838709
$ws.Range("F21").Value = 190
$ws.Range("G21").Value = 152
$ws.Range("H21").Value = 25
$ws.Range("I21").Value = 2313
$ws.Range("J21").Value = 2013
$ws.Range("K21").Value = 14.903129657228
$ws.Range("L21").Value = 23.888591322978
$ws.Range("M21").Value = 15.017404276479
$ws.Range("N21").Value = -65.487914055505
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 47
$ws.Range("J22").Value = 45
$ws.Range("K22").Value = 4.444444444444
$ws.Range("L22").Value = 51.612903225806
$ws.Range("M22").Value = 88
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 65
$ws.Range("E24").Value = -41.538461538461
$ws.Range("F24").Value = 179
$ws.Range("G24").Value = 137
$ws.Range("H24").Value = 30.656934306569
$ws.Range("I24").Value = 2373
$ws.Range("J24").Value = 2466
$ws.Range("K24").Value = -3.771289537712
$ws.Range("L24").Value = 66.760365425158
$ws.Range("M24").Value = 99.244332493702
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -35.714285714285
$ws.Range("F25").Value = 59
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 766
$ws.Range("J25").Value = 730
$ws.Range("K25").Value = 4.931506849315
$ws.Range("L25").Value = 5.075445816186
$ws.Range("M25").Value = -7.710843373493
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 59
$ws.Range("J26").Value = 68
$ws.Range("K26").Value = -13.235294117647
$ws.Range("L26").Value = 43.90243902439
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 400
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 28.571428571428
$ws.Range("I27").Value = 96
$ws.Range("J27").Value = 82
$ws.Range("K27").Value = 17.073170731707
$ws.Range("L27").Value = -4.950495049504
$ws.Range("G28").Value = 1
$ws.Range("M28").Value = -60.526315789473
$ws.Range("G29").Value = 1
$ws.Range("M29").Value = -60
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = 25
$ws.Range("L30").Value = 66.666666666666
